{"js": "// The document contains 4 paragraphs with the old Slovenian \"2018: Datumi\n// kampanje za opazovanje Perseus: 30. oktobra - 8. novembra in 29. novembra\n// - 8. decembra\" text (some split across multiple runs, one followed by\n// extra \"Obvezno izpolnite...\" runs). Each such paragraph must end up with\n// a single, plain (unformatted) run containing the new text.\nconst NEW_TEXT =\n  \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Bootes: 14. in 23. maja, 13. in 22. junij, 12. in 21. julija\";\nconst OLD_MARKER = \"Datumi kampanje\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Collect indices first (mutating while iterating the same loaded collection\n// is safe since we only read .text, already loaded).\nconst targetIndexes = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(OLD_MARKER) !== -1) {\n    targetIndexes.push(i);\n  }\n}\n\nfor (const idx of targetIndexes) {\n  const para = paragraphs.items[idx];\n  // Remove every run (and their formatting) from the paragraph...\n  para.clear();\n  // ...then insert fresh, unformatted text as the paragraph's sole run.\n  para.insertText(NEW_TEXT, Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# The document contains 4 paragraphs with the old Slovenian \"2018: Datumi\n# kampanje za opazovanje Perseus: 30. oktobra - 8. novembra in 29. novembra\n# - 8. decembra\" text (some split across multiple runs, one followed by\n# extra \"Obvezno izpolnite...\" runs). Each such paragraph must end up with\n# a single, plain (unformatted) run containing the new text.\n$d = $word.ActiveDocument\n$NEW_TEXT = \"m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Bootes: 14. in 23. maja, 13. in 22. junij, 12. in 21. julija\"\n$OLD_MARKER = \"Datumi kampanje\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    if ($t -like \"*$OLD_MARKER*\") {\n        # Delete all paragraph content (everything up to, but excluding,\n        # the trailing paragraph mark) so no run/character formatting\n        # survives for the new text.\n        $contentRange = $d.Range($p.Range.Start, $p.Range.End - 1)\n        $contentRange.Delete()\n\n        # Re-fetch the (now empty) paragraph's range and insert the\n        # replacement text as a brand-new, unformatted run.\n        $p2 = $d.Paragraphs($i)\n        $insertRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)\n        $insertRange.InsertAfter($NEW_TEXT)\n    }\n}\n"}
